$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 67.8255122017956
$ws.Range("C2").Value = 66.435391002076
$ws.Range("D2").Value = 69.2156334015152
$ws.Range("B3").Value = 58.2213755593942
$ws.Range("C3").Value = 51.6233686680771
$ws.Range("D3").Value = 64.8193824507112
$ws.Range("C4").Value = 59.3069211203831
$ws.Range("D4").Value = 69.608124611224
$ws.Range("C5").Value = 76.4362667906174
$ws.Range("D5").Value = 83.6206810157816
$ws.Range("C6").Value = 59.610099416897
$ws.Range("D6").Value = 68.899933990014
$ws.Range("C7").Value = 41.3914102131016
$ws.Range("D7").Value = 52.5070576892007
$ws.Range("C8").Value = 61.4179413522406
$ws.Range("D8").Value = 75.8556927259095
$ws.Range("C9").Value = 60.6944421035916
$ws.Range("D9").Value = 77.9485000207496
$ws.Range("C10").Value = 64.052181693643
$ws.Range("D10").Value = 76.1968857597093
$ws.Range("C11").Value = 82.6459289308268
$ws.Range("D11").Value = 92.585982284678
$ws.Range("C12").Value = 84.7078721493097
$ws.Range("D12").Value = 97.1415898437397
$ws.Range("C13").Value = 59.713866721804
$ws.Range("D13").Value = 67.5160599020405
$ws.Range("B14").Value = 78.5390250556913
$ws.Range("C14").Value = 75.1487917601202
$ws.Range("D14").Value = 81.9292583512623
$ws.Range("C15").Value = 45.1584513120294
$ws.Range("D15").Value = 54.2327991495909
$ws.Range("B16").Value = 55.3746767090015
$ws.Range("C16").Value = 54.4859318721968
$ws.Range("D16").Value = 56.2634215458062
$ws.Range("B17").Value = 45.2808231861887
$ws.Range("C17").Value = 41.2487754709432
$ws.Range("D17").Value = 49.3128709014343
$ws.Range("B18").Value = 60.8336262307661
$ws.Range("C18").Value = 57.7966356298784
$ws.Range("D18").Value = 63.8706168316539
$ws.Range("C19").Value = 75.2134161195928
$ws.Range("D19").Value = 79.6398147110472
$ws.Range("C20").Value = 57.1778279086813
$ws.Range("D20").Value = 62.7218865374742
$ws.Range("C21").Value = 42.6035610321692
$ws.Range("D21").Value = 49.0492070796423
$ws.Range("C22").Value = 43.7996249188381
$ws.Range("D22").Value = 53.0646946537023
$ws.Range("C23").Value = 56.2066845768809
$ws.Range("D23").Value = 66.2156236688857
$ws.Range("C24").Value = 53.5489195492056
$ws.Range("D24").Value = 60.9173713018241
$ws.Range("C25").Value = 80.8499576095317
$ws.Range("D25").Value = 88.4357636282268
$ws.Range("C26").Value = 79.959942270549
$ws.Range("D26").Value = 90.3559851914384
$ws.Range("C27").Value = 51.2958490153824
$ws.Range("D27").Value = 55.9624280812738
$ws.Range("C28").Value = 55.5458334564236
$ws.Range("D28").Value = 60.2907784024995
$ws.Range("C29").Value = 32.7455137903951
$ws.Range("D29").Value = 37.6780054834294
$ws.Range("B30").Value = 64.7210522905015
$ws.Range("C30").Value = 63.0016987755704
$ws.Range("D30").Value = 66.4404058054326
$ws.Range("B31").Value = 53.6064814821634
$ws.Range("C31").Value = 45.1469432318865
$ws.Range("D31").Value = 62.0660197324404
$ws.Range("C32").Value = 58.8617915566398
$ws.Range("D32").Value = 70.5234969493626
$ws.Range("C33").Value = 76.9119737744376
$ws.Range("D33").Value = 85.0996864571568
$ws.Range("C34").Value = 62.7478079397744
$ws.Range("D34").Value = 73.0553545713773
$ws.Range("C35").Value = 49.8921441734352
$ws.Range("D35").Value = 62.6023719162634
$ws.Range("C36").Value = 54.6119886877826
$ws.Range("D36").Value = 71.8688303029242
$ws.Range("C37").Value = 62.9705398426489
$ws.Range("D37").Value = 81.6145802716471
$ws.Range("C38").Value = 60.8508236706679
$ws.Range("D38").Value = 75.1381165881873
$ws.Range("C39").Value = 80.0840316706184
$ws.Range("D39").Value = 94.3597930822516
$ws.Range("C40").Value = 81.4193306110471
$ws.Range("D40").Value = 98.8293312295551
$ws.Range("C41").Value = 59.2211594131443
$ws.Range("D41").Value = 68.2353311225131
$ws.Range("C42").Value = 69.6919360272835
$ws.Range("D42").Value = 78.1348070970025
$ws.Range("C43").Value = 37.4998334468133
$ws.Range("D43").Value = 47.9212342432281
$ws.Range("B44").Value = 67.1354010141054
$ws.Range("C44").Value = 65.427782749948
$ws.Range("D44").Value = 68.8430192782629
$ws.Range("B45").Value = 61.1177952962955
$ws.Range("C45").Value = 53.2562585890599
$ws.Range("D45").Value = 68.9793320035311
$ws.Range("C46").Value = 58.6159445104138
$ws.Range("D46").Value = 70.5658444950856
$ws.Range("C47").Value = 86.0077669701419
$ws.Range("D47").Value = 92.4144253246908
$ws.Range("C48").Value = 66.4627623231789
$ws.Range("D48").Value = 76.9906166015781
$ws.Range("C49").Value = 50.0835950826615
$ws.Range("D49").Value = 63.7547446715367
$ws.Range("C50").Value = 50.7223889638352
$ws.Range("D50").Value = 69.4799371741746
$ws.Range("C51").Value = 65.3663760773419
$ws.Range("D51").Value = 83.8371788933479
$ws.Range("B52").Value = 68.6349696406961
$ws.Range("C52").Value = 61.8659844146273
$ws.Range("D52").Value = 75.4039548667649
$ws.Range("B53").Value = 86.9268829131096
$ws.Range("C53").Value = 79.9246344114365
$ws.Range("D53").Value = 93.9291314147826
$ws.Range("C54").Value = 82.8187124167857
$ws.Range("D54").Value = 99.0061195585191
$ws.Range("C55").Value = 52.0736629975913
$ws.Range("D55").Value = 61.4699106202318
$ws.Range("C56").Value = 72.3001749304455
$ws.Range("D56").Value = 80.7144135175635
$ws.Range("C57").Value = 34.1895696793302
$ws.Range("D57").Value = 44.9682067069441
$ws.Range("B58").Value = 68.4899159160604
$ws.Range("C58").Value = 66.8593150984853
$ws.Range("D58").Value = 70.1205167336354
$ws.Range("B59").Value = 63.1815934774769
$ws.Range("C59").Value = 55.6156662398226
$ws.Range("D59").Value = 70.7475207151311
$ws.Range("C60").Value = 64.7385781207248
$ws.Range("D60").Value = 76.0718194276089
$ws.Range("C61").Value = 84.2175744593211
$ws.Range("D61").Value = 91.2525943997292
$ws.Range("C62").Value = 69.2133129270076
$ws.Range("D62").Value = 79.1962245022285
$ws.Range("C63").Value = 50.5308726568245
$ws.Range("D63").Value = 62.6158997105451
$ws.Range("B64").Value = 65.2899469820092
$ws.Range("C64").Value = 56.3334248259056
$ws.Range("D64").Value = 74.2464691381127
$ws.Range("C65").Value = 76.1349441980019
$ws.Range("D65").Value = 91.2844966347573
$ws.Range("C66").Value = 62.6694697105599
$ws.Range("D66").Value = 76.4612164325803
$ws.Range("C67").Value = 84.0298059736881
$ws.Range("D67").Value = 96.4702275280408
$ws.Range("C68").Value = 85.5750992589754
$ws.Range("D68").Value = 100.829035227341
$ws.Range("C69").Value = 63.6547257655049
$ws.Range("D69").Value = 72.1893534796417
$ws.Range("C70").Value = 76.40218349686
$ws.Range("D70").Value = 83.8782787745212
$ws.Range("C71").Value = 41.0271318868865
$ws.Range("D71").Value = 50.798946287091
$ws.Range("B72").Value = 61.8798692282585
$ws.Range("C72").Value = 60.1485598973454
$ws.Range("D72").Value = 63.6111785591716
$ws.Range("B73").Value = 52.8533714542478
$ws.Range("C73").Value = 44.7007246364643
$ws.Range("D73").Value = 61.0060182720313
$ws.Range("C74").Value = 61.1969262822731
$ws.Range("D74").Value = 73.1694027676559
$ws.Range("B75").Value = 80.2386934272549
$ws.Range("C75").Value = 75.8815557174006
$ws.Range("D75").Value = 84.5958311371091
$ws.Range("C76").Value = 64.7000509080305
$ws.Range("D76").Value = 74.9824319234582
$ws.Range("C77").Value = 41.7907434576787
$ws.Range("D77").Value = 54.9733615794522
$ws.Range("C78").Value = 47.1340090898982
$ws.Range("D78").Value = 65.6064879566911
$ws.Range("C79").Value = 59.5431942466013
$ws.Range("D79").Value = 78.017354569529
$ws.Range("C80").Value = 62.0095762398281
$ws.Range("D80").Value = 75.7856419316606
$ws.Range("C81").Value = 79.5477714514187
$ws.Range("D81").Value = 93.9284952012993
$ws.Range("C82").Value = 81.62619155425
$ws.Range("D82").Value = 99.15323888893
$ws.Range("C83").Value = 57.2885087709571
$ws.Range("D83").Value = 66.4995937093343
$ws.Range("C84").Value = 65.4913298875108
$ws.Range("D84").Value = 74.3042387074987
$ws.Range("C85").Value = 36.2604272201354
$ws.Range("D85").Value = 46.2283594196645
